$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) tc_01 (sheet1): remove the 3rd (last) row that held the stray
#    "q34234"/"dfsdf" test values, and leave the selection positioned
#    on the now-empty row 3 as the next entry point.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("tc_01")
$null = $ws1.Rows.Item(3).Delete()
$null = $ws1.Range("A3:B3").Select()

# ---------------------------------------------------------------------
# 2) Add a brand new "testData" worksheet after the last existing sheet.
#    Duplicate an existing (uncustomized) sheet and wipe it so the new
#    sheet doesn't pick up any unwanted default formatting metadata,
#    then rename + clear it.
# ---------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("tc_06")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$templateSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$ws7 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7.Name = "testData"
$ws7.Cells.Clear()

# ---------------------------------------------------------------------
# 3) Populate the testData sheet with the new ZOHO.com test content.
#    Values are written in the same order they were first introduced
#    so the shared-string table ends up ordered the same way.
# ---------------------------------------------------------------------
$ws7.Range("A1").Value = "tc_02"

$ws7.Range("A2").Value = "username"
$ws7.Range("B2").Value = "password"

$ws7.Range("A3").Value = 8744954505
$ws7.Range("B3").Value = 12123

$ws7.Range("A5").Value = "tc_01"

$ws7.Range("A6").Value = "username"
$ws7.Range("B6").Value = "password"

$ws7.Range("A7").Value = "esdf"
$ws7.Range("B7").Value = "sdwerwe"

$ws7.Range("A8").Value = "sdfsdf"
$ws7.Range("B8").Value = "qwewer"

$ws7.Range("A9").Value = "fsfsd"
$ws7.Range("B9").Value = "dsdfsdf"

$ws7.Range("C6").Value = "email"
$ws7.Range("D6").Value = "firstname"
$ws7.Range("E6").Value = "lastaname"

# Column A was sized to fit its contents.
$ws7.Columns.Item(1).ColumnWidth = 10.1666666667

# Final selection/active cell on the new sheet.
$null = $ws7.Range("E6").Select()

Write-Host "Edit complete"
